# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Marzo de 2020 a las 01:59"

# 2. Update Estados Unidos (row 4) stats
$ws.Range("B4").Value = 103729
$ws.Range("C4").Value = 18294
$ws.Range("D4").Value = 2522
$ws.Range("E4").Value = 99514
$ws.Range("F4").Value = 2463
$ws.Range("G4").Value = 398
$ws.Range("H4").Value = 1693

# 3. Update Australia (row 22) stats
$ws.Range("B22").Value = 3378
$ws.Range("C22").Value = 328
$ws.Range("D22").Value = 170
$ws.Range("E22").Value = 3195
$ws.Range("F22").Value = 23

# 4. Reorder Panama ahead of Singapur/Crucero (rows 46-48) with updated data
# Row 46: was Singapur(732,49,183,547,17,0,2) -> becomes Panama with new data
$ws.Range("A46").Value = "Panama"
$ws.Range("B46").Value = 786
$ws.Range("C46").Value = 112
$ws.Range("D46").Value = 2
$ws.Range("E46").Value = 770
$ws.Range("F46").Value = 20
$ws.Range("G46").Value = 5
$ws.Range("H46").Value = 14

# Row 47: was Crucero(712,0,597,105,15,0,10) -> becomes Singapur with old row46 data
$ws.Range("A47").Value = "Singapur"
$ws.Range("B47").Value = 732
$ws.Range("C47").Value = 49
$ws.Range("D47").Value = 183
$ws.Range("E47").Value = 547
$ws.Range("F47").Value = 17
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 2

# Row 48: was Panama(674,0,2,663,20,0,9) -> becomes Crucero with old row47 data
$ws.Range("A48").Value = "Crucero"
$ws.Range("B48").Value = 712
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 597
$ws.Range("E48").Value = 105
$ws.Range("F48").Value = 15
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 10

# Row 49 (Peru) is unchanged.

# 5. Update Honduras (row 115) new cases death column
$ws.Range("F115").Value = 4
